$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Adiciona o projeto do módulo 'Construindo páginas para internet com bootstrap'
# (linha 10) e o respectivo certificado (linha 11), preenchendo as horas.
$ws.Range("D9").Copy()
$ws.Range("D10:D11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("D10").Value = 6
$ws.Range("D11").Value = 4

# Atualiza a seleção ativa para refletir a última célula editada
$ws.Range("D12").Select()

$wb.Save()
